$wb = $excel.ActiveWorkbook

# Remember which sheet is active (SDQ) so we can restore it at the end -
# the edits below touch "Service Contacts" and "Practitioners" without
# the SDQ sheet losing its place as the active/selected tab.
$wsActive = $wb.ActiveSheet

# --- Service Contacts sheet: widen column A a touch + move the selection ---
$wsSC = $wb.Worksheets.Item("Service Contacts")
$wsSC.Activate()
# ColumnWidth = 13.6667 round-trips to the saved "width" of 14.5.
$wsSC.Columns.Item(1).ColumnWidth = 13.666666666666666
$wsSC.Range("D3").Select()

# --- Practitioners sheet: new data row + column widths + selection ---
$wsP = $wb.Worksheets.Item("Practitioners")
$wsP.Activate()
# ColumnWidth values chosen so the saved "width" lands on/near 14.6640625,
# 13 and 12.83203125 respectively (A, C, F).
$wsP.Columns.Item(1).ColumnWidth = 13.833333333333332
$wsP.Columns.Item(3).ColumnWidth = 12.166666666666666
$wsP.Columns.Item(6).ColumnWidth = 12.0

$wsP.Range("A6").Value = "PHN999:NFP02"
$wsP.Range("B6").Value = "P01"
$wsP.Range("C6").Value = 8
$wsP.Range("D6").Value = 1
$wsP.Range("E6").Value = 1973
$wsP.Range("F6").Value = 2
$wsP.Range("G6").Value = 1
$wsP.Range("H6").Value = 1
$wsP.Range("I6").Value = "tag1"

$wsP.Columns.Item(7).Select()

# Restore the originally-active sheet (SDQ) so the workbook-level active
# tab / this sheet's tabSelected flag are unchanged.
$wsActive.Activate()

$wb.Save()
